$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the marking scheme (right-answer mark) and the recomputed totals
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 65
$ws.Range("E12").Value = "65/140"
